# "Improved names and tests"
#
# Updates two of the numeric-looking text values in the test-data sheet and
# moves the active-cell selection from C4 to E4.
#
#   E3: "2.00001" -> "2.00501"
#   E4: "3.00001" -> "3.00401"
#   Selection: C4 -> E4

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = "2.00501"
$ws.Range("E4").Value = "3.00401"

$ws.Range("E4").Select()
